$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = "2026-01-16 06:40:14"
}
